# "journal de bord" - ajout de page de deconnexion si besoin, probablement enlever si innutile
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New column C content --------------------------------------------------
# Write these two new strings first so the shared-strings table ends up in
# the same append order as the target workbook (revision/connexion before
# the corrected "comprention" string below).
$ws.Range("C2").Value = "révision de la BDD & ajout de quelque valeurs"
$ws.Range("C3").Value = "connexion a la BDD faites dans l'index"

# Fix the typo in B4 ("compreention" -> "comprention"); this adds a brand
# new shared string (appended last), matching the target ordering.
$ws.Range("B4").Value = "commencer comprention matos"

# C1: second date, one day after B1 (2021-10-04 -> 2021-10-05), formatted
# the same way as B1. Copy/paste the format from B1 so the existing date
# number format style gets reused instead of creating a new one.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C1").Value = 44474

# New column width for column C.
$ws.Columns.Item(3).ColumnWidth = 41.25

# Match the saved selection shown in the target file.
$ws.Range("C7").Select()

$wb.Save()
